# cryptos.xlsx refresh -- GitHub Actions scheduled price/volume update.
# Rows 40-42 also get re-ranked (coin name/link/price/volume all move).
#
# A leading apostrophe forces Excel to keep a numeric-looking string (e.g.
# "551.40" or "0.0000170") as literal text instead of silently coercing it
# to a Number (which would drop trailing zeros or flip to scientific notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '67.470.34'
$ws.Range('E2').Value = '  -2.31%  '
# Row 3: Ethereum
$ws.Range('D3').Value = '2.423.80'
$ws.Range('E3').Value = '  -1.85%  '
# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.01%  '
# Row 5: BNB
$ws.Range('D5').Value = '''551.40'
$ws.Range('E5').Value = '  -1.29%  '
# Row 6: Solana
$ws.Range('D6').Value = '''158.53'
$ws.Range('E6').Value = '  -2.70%  '
# Row 7: USDC
$ws.Range('E7').Value = '  +0.03%  '
# Row 8: XRP
$ws.Range('D8').Value = '''0.506'
$ws.Range('E8').Value = '  +0.24%  '
# Row 9: Dogecoin
$ws.Range('D9').Value = '''0.158'
$ws.Range('E9').Value = '  +4.30%  '
# Row 10: TRON
$ws.Range('E10').Value = '  -1.13%  '
# Row 11: Cardano
$ws.Range('D11').Value = '''0.328'
$ws.Range('E11').Value = '  -2.07%  '
# Row 12: Toncoin
$ws.Range('D12').Value = '''4.77'
$ws.Range('E12').Value = '  -1.24%  '
# Row 13: WrappedBTC
$ws.Range('D13').Value = '67.402.73'
$ws.Range('E13').Value = '  -2.32%  '
# Row 14: ShibaInu
$ws.Range('D14').Value = '''0.0000170'
$ws.Range('E14').Value = '  +0.01%  '
# Row 15: Avalanche
$ws.Range('D15').Value = '''22.91'
$ws.Range('E15').Value = '  -3.00%  '
# Row 16: Chainlink
$ws.Range('D16').Value = '''10.35'
$ws.Range('E16').Value = '  -3.87%  '
# Row 17: BitcoinCash
$ws.Range('D17').Value = '''328.26'
$ws.Range('E17').Value = '  -4.20%  '
# Row 18: Uniswap
$ws.Range('D18').Value = '''6.82'
$ws.Range('E18').Value = '  -3.24%  '
# Row 19: Polkadot
$ws.Range('D19').Value = '''3.78'
$ws.Range('E19').Value = '  -0.78%  '
# Row 20: Dai
$ws.Range('E20').Value = '  -0.52%  '
# Row 21: SuiNetwork
$ws.Range('D21').Value = '''1.84'
$ws.Range('E21').Value = '  -1.97%  '
# Row 22: Litecoin
$ws.Range('D22').Value = '''65.72'
$ws.Range('E22').Value = '  -1.94%  '
# Row 23: NEARProtocol
$ws.Range('D23').Value = '''3.60'
$ws.Range('E23').Value = '  -2.33%  '
# Row 24: Aptos
$ws.Range('D24').Value = '''8.01'
$ws.Range('E24').Value = '  -1.79%  '
# Row 25: PEPE
$ws.Range('D25').Value = '0.0₃0801'
$ws.Range('E25').Value = '  -2.37%  '
# Row 26: InternetComputer(DFINITY)
$ws.Range('D26').Value = '''7.01'
$ws.Range('E26').Value = '  -2.58%  '
# Row 27: FirstDigitalUSD
$ws.Range('E27').Value = '  +0.04%  '
# Row 28: Bittensor
$ws.Range('D28').Value = '''415.44'
$ws.Range('E28').Value = '  -5.22%  '
# Row 29: Fetch.AI
$ws.Range('D29').Value = '''1.12'
$ws.Range('E29').Value = '  -2.15%  '
# Row 30: PancakeSwap
$ws.Range('D30').Value = '''1.59'
$ws.Range('E30').Value = '  -1.69%  '
# Row 31: Monero
$ws.Range('D31').Value = '''159.33'
$ws.Range('E31').Value = '  +1.29%  '
# Row 32: WhiteBITCoin
$ws.Range('E32').Value = '  -0.77%  '
# Row 33: USDe
$ws.Range('E33').Value = '  -0.02%  '
# Row 34: EthereumClassic
$ws.Range('D34').Value = '''17.74'
$ws.Range('E34').Value = '  -0.93%  '
# Row 35: Kaspa
$ws.Range('E35').Value = '  -4.49%  '
# Row 36: PolygonEcosystemToken
$ws.Range('D36').Value = '''0.293'
$ws.Range('E36').Value = '  -3.56%  '
# Row 37: RenderToken
$ws.Range('D37').Value = '''4.22'
$ws.Range('E37').Value = '  -5.41%  '
# Row 38: Stacks
$ws.Range('D38').Value = '''1.45'
$ws.Range('E38').Value = '  -1.84%  '
# Row 39: ImmutableX
$ws.Range('D39').Value = '''1.07'
$ws.Range('E39').Value = '  -3.18%  '
# Row 40: dogwifhat (was Aave)
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = '''1.98'
$ws.Range('E40').Value = '  -4.95%  '
# Row 41: Filecoin (was dogwifhat)
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '''3.30'
$ws.Range('E41').Value = '  -1.62%  '
# Row 42: Aave (was Filecoin)
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '''129.41'
$ws.Range('E42').Value = '  -2.99%  '
# Row 43: Cronos
$ws.Range('D43').Value = '''0.0707'
$ws.Range('E43').Value = '  -1.61%  '
# Row 44: ARBITRUM
$ws.Range('D44').Value = '''0.475'
$ws.Range('E44').Value = '  -1.96%  '
# Row 45: Mantle
$ws.Range('D45').Value = '''0.553'
$ws.Range('E45').Value = '  -2.07%  '
# Row 46: Stellar
$ws.Range('D46').Value = '''0.0912'
$ws.Range('E46').Value = '  +0.56%  '
# Row 47: BitgetToken
$ws.Range('E47').Value = '  -0.26%  '
# Row 48: Optimism
$ws.Range('D48').Value = '''1.33'
$ws.Range('E48').Value = '  -8.09%  '
# Row 49: InjectiveProtocol
$ws.Range('D49').Value = '''16.48'
$ws.Range('E49').Value = '  -2.70%  '
# Row 50: BabyDogeCoin
$ws.Range('D50').Value = '0.0₆0203'
$ws.Range('E50').Value = '  +0.04%  '
# Row 51: Hedera
$ws.Range('D51').Value = '''0.0426'
$ws.Range('E51').Value = '  -1.03%  '
